$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the existing header formatting
# used by B1:G1 (bold font, thin box border, centered/top aligned).
$h1 = $ws.Cells.Item(1, 8)
$h1.Value = "Save"

$h1.Font.Bold = $true
$h1.HorizontalAlignment = -4108   # xlCenter
$h1.VerticalAlignment = -4160     # xlTop

$h1.Borders.LineStyle = 1         # xlContinuous
$h1.Borders.Weight = 2            # xlThin

# Fill in the new "Save" column values for the data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
